# adding averages and more checks
#
# Updates the Training Dashboard's "PERIOD TO EXPIRE" / "LAST UPDATE" columns
# to reflect a refreshed check date (08-Sep-2025 -> 16-Sep-2025, an 8-day
# shift, with PERIOD TO EXPIRE reduced by 8 days to match), refreshes the
# Exam Dashboard's comments column with a more descriptive message (and
# widens that column to fit), and restyles the title / header bands to use
# bold white text (instead of a big title font + separate plain bold header
# font).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Training Dashboard: refresh "LAST UPDATE" date and recompute the
#    "PERIOD TO EXPIRE" day-count (each dropped by 8, matching the new
#    8-day-later check date) for every data row (3-14).
# ---------------------------------------------------------------------
$training = $wb.Worksheets.Item("Training Dashboard")

$periodToExpire = @{
    3  = 213
    4  = 359
    5  = 404
    6  = 404
    7  = 210
    8  = 353
    9  = 113
    10 = -343
    11 = 679
    12 = 679
    13 = 679
    14 = 679
}

foreach ($row in 3..14) {
    # Column H = PERIOD TO EXPIRE (plain numeric, keeps existing style).
    $training.Cells.Item($row, 8).Value = $periodToExpire[$row]

    # Column I = LAST UPDATE. The sheet stores this as literal text (not a
    # real Excel date), so force text interpretation before writing, to
    # avoid Excel's automatic "looks like a date" conversion.
    $cell = $training.Cells.Item($row, 9)
    $cell.NumberFormat = "@"
    $cell.Value = "16-Sep-2025"
}

# ---------------------------------------------------------------------
# 2. Exam Dashboard: friendlier comment text + widen column E to fit it.
# ---------------------------------------------------------------------
$exam = $wb.Worksheets.Item("Exam Dashboard")

$exam.Range("E3").Value = "date is valid"
$exam.Range("E4").Value = "date is valid"

# Raw OOXML column <col width="..."> runs ~0.8333 wider than the COM
# ColumnWidth figure (character-width vs. pixel-padded units), so back that
# offset out to land exactly on width=15.
$exam.Columns("E").ColumnWidth = 15 - (5 / 6)

# ---------------------------------------------------------------------
# 3. Restyle the title banner and column-header band on both sheets to
#    bold white text (dropping the separate oversized 14pt title font and
#    the plain-bold header font in favour of one shared bold/white font).
# ---------------------------------------------------------------------
foreach ($ws in @($training, $exam)) {
    $lastCol = $ws.UsedRange.Columns.Count  # header band spans the full used width, including trailing blank-but-styled cells

    $title = $ws.Range("A1")
    $title.Font.Bold = $true
    $title.Font.Size = 11
    $title.Font.Color = 16777215  # white

    $header = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item(2, $lastCol))
    $header.Font.Bold = $true
    $header.Font.Color = 16777215  # white
}

Write-Output "applied training/exam dashboard refresh"
